$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update Tiempo_Mínimo (E), Tiempo_Máximo (F), Tiempo_Promedio (G) values
# Row 2
$ws.Range("E2").Value = 0.000223218
$ws.Range("F2").Value = 0.017757168
$ws.Range("G2").Value = 0.00036514324190000004

# Row 3
$ws.Range("E3").Value = 0.002460208
$ws.Range("F3").Value = 0.006888562
$ws.Range("G3").Value = 0.0028385706895963617

# Row 4
$ws.Range("E4").Value = 0.010465599
$ws.Range("F4").Value = 0.015279491
$ws.Range("G4").Value = 0.011949152739234449
